# assignmenttable.xlsx update:
# - Convert the "start date" (C) / "due date" (D) columns from text labels
#   (e.g. "06.26(수)") to real Excel dates, formatted with a custom
#   m"/"d;@ number format.
# - Re-point the assignment title / content / related-lecture columns
#   (E, F, G) at the correct rows (the data had drifted out of alignment).
# - Update the sheet view's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the custom date number format to the whole start/due-date block.
$ws.Range("C2:D8").NumberFormat = 'm"/"d;@'

# Row 2 - [플러스엑스]UX 디자인 컨셉과 전략
$ws.Range("C2").Value = 45463
$ws.Range("D2").Value = 45469
$ws.Range("E2").Value = "전략문서 작성하기"
$ws.Range("F2").Value = "본인이 하고 싶은 서비스 선정 후 전략 문서 제작하기`n목차`n1. 서비스 소개 및 분석`n2. 데스크 리서치`n3. 사용자 케이스 설정`n4. 벤치마킹 서비스 분석`n4. 스토리 라인 설정`n`n장표를 PDF 로 제출(장표 분량 자유)"
$ws.Range("G2").Value = "CH03-01 UX 전략 문서`nCH03-02 UX 전략 문서를 만드는 과정"

# Row 3 - 모두를 위한 피그마 119개 실습으로 완전 정복
$ws.Range("C3").Value = 45470
$ws.Range("D3").Value = 45473
$ws.Range("E3").Value = "반응형 페이지 만들기"
$ws.Range("F3").Value = "Figma 강의를 듣고 반응형 페이지 제작하기 [Autolayout]"
$ws.Range("G3").Value = "Ch09-09. [연습] 반응형 페이지 만들기"

# Row 4 - UX 유저 리서치, 사용성 평가, UX 데이터 모델링
$ws.Range("C4").Value = 45481
$ws.Range("D4").Value = 45499
$ws.Range("E4").Value = "UX 그룹 과제 최종 결과물"
$ws.Range("F4").Value = "[그룹 과제]`n`n■ 1~2회 차`n1. APPROACH 수립, 가설 수립`n2. 사용자 조사 기초 살계`n`n■ 3회 차`n1. 조사 대상자 모집`n2. 조사 질문지 작성`n`n■ 4회 차`n1. 조사 수행 및 결과 정리`n`n■ 5회 차`n1. UX 패턴/UX 모델링`n2. UX 전략`n`n■ 6회 차`n1. 발표회`n`n위 진행에 따른 최종 결과물 제출"
$ws.Range("G4").Value = "Ch 01. 일러스트 입문`n04. 회전 툴을 이용한 톱니바퀴 그리기"

# Row 5 - UI Design
$ws.Range("C5").Value = 45482
$ws.Range("D5").Value = 45485
$ws.Range("E5").Value = "UI 온라인 강의 1차 과제"
$ws.Range("F5").Value = "[UI 디자인 기초 지식]`n1. 최근 자신이 작업했던 경험을 비추어 볼 때, 관련 지식들이 없어서 힘들었던 점을 최대한 생각해 내어 적어보도록 하자.`n 그리고 이런 상황에서 어떻게 그때 그때 넘겼는지도 함께 적어보자."
$ws.Range("G5").Value = "CH01-02. UI 디자인 기초 지식"

# Row 6 - 초격차 패키지 : 한 번에 끝내는 디자인 툴
$ws.Range("C6").Value = 45489
$ws.Range("D6").Value = 45491
$ws.Range("E6").Value = "포토샵 온라인 강의 2차 과제"
$ws.Range("F6").Value = "Layer Style을 활용한 상품 광고 이미지 제작하기"
$ws.Range("G6").Value = "06. Layer Style을 활용한 이미지 제작 - 글자 윤곽선"

# Row 7 - 초격차 패키지 : 한 번에 끝내는 디자인 툴 (trailing space variant)
$ws.Range("C7").Value = 45496
$ws.Range("D7").Value = 45498
$ws.Range("E7").Value = "UI 온라인 강의 3차 과제"
$ws.Range("F7").Value = "Pathfinder 패널을 이용한 수박 그리기"
$ws.Range("G7").Value = "22. 채색하는 방법을 응용한 화장품 그리기"

# Row 8 - 초격차 패키지 : 한 번에 끝내는 디자인 툴
$ws.Range("C8").Value = 45503
$ws.Range("D8").Value = 45508
$ws.Range("E8").Value = "UI 온라인 강의 4차 과제"
$ws.Range("F8").Value = "다양한 기능을 활용한 배경 그리기"
$ws.Range("G8").Value = "24. 다양한 기능을 활용한 배경 그리기_1"

# Move the active selection to E4 (matches the saved view state in the file).
$ws.Range("E4").Select()
